$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# serialVersion row: it is now managed/controlled by CN (Coordinating Node)
# replication processes rather than an MN-service-subject.
$ws.Range("C2").Value = "CN"
$ws.Range("F2").Value = "CN replication processes"

# Move the active selection/cursor to E7 (matches the saved cursor position).
$ws.Range("E7").Select()
